$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 188, pushing the existing rows
# 188..220 down to 190..222 (formatting/styles carried from row above).
$ws.Rows.Item(188).Insert()
$ws.Rows.Item(189).Insert()

# Row 188: new "Primera" quality record for Fecha 2021-11-05 (44505)
$ws.Cells.Item(188, 1).Value = 4
$ws.Cells.Item(188, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(188, 3).Value = "Los Lagos"
$ws.Cells.Item(188, 4).Value = 44505
$ws.Cells.Item(188, 5).Value = 10
$ws.Cells.Item(188, 6).Value = 100112023
$ws.Cells.Item(188, 7).Value = "Brócoli"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 500
$ws.Cells.Item(188, 11).Value = 1200
$ws.Cells.Item(188, 12).Value = 1200
$ws.Cells.Item(188, 13).Value = 1200
$ws.Cells.Item(188, 14).Value = "$/unidad"
$ws.Cells.Item(188, 15).Value = "Región Metropolitana"
$ws.Cells.Item(188, 16).Value = 1200
$ws.Cells.Item(188, 17).Value = 1
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Row 189: new "Segunda" quality record for the same Fecha 2021-11-05 (44505)
$ws.Cells.Item(189, 1).Value = 4
$ws.Cells.Item(189, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(189, 3).Value = "Los Lagos"
$ws.Cells.Item(189, 4).Value = 44505
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 100112023
$ws.Cells.Item(189, 7).Value = "Brócoli"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Segunda"
$ws.Cells.Item(189, 10).Value = 500
$ws.Cells.Item(189, 11).Value = 1000
$ws.Cells.Item(189, 12).Value = 1000
$ws.Cells.Item(189, 13).Value = 1000
$ws.Cells.Item(189, 14).Value = "$/unidad"
$ws.Cells.Item(189, 15).Value = "Región Metropolitana"
$ws.Cells.Item(189, 16).Value = 1000
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"
